$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the floating point rounding on A44's timestamp value
$ws.Range("A44").Value = 45055.04609740741

# New chat log rows to append
$rows = @(
    @{ A = 45056.10041508102; B = "Hi"; C = "Howdy! How can I help you today?" },
    @{ A = 45056.10045188657; B = "Can you help me?"; C = "Can you please clarify your question?" },
    @{ A = 45056.10048769676; B = "my name is matthew"; C = "Hello Matthew. How can I help you?" },
    @{ A = 45056.10053506945; B = "How can you help me?"; C = "Can you please clarify your question?" },
    @{ A = 45056.10056601852; B = "Can you help me?"; C = "Can you please clarify your question?" },
    @{ A = 45056.10105740741; B = "Hi"; C = "Howdy! How can I help you today?" },
    @{ A = 45056.10109155093; B = "Can you help me?"; C = "Choosing the right ransomware solution depends on a number of factors, including the size of your organization, your security needs, and your budget. I can help you narrow down your options and choose the solution that is right for you. Would you like me to provide you with more information about our products?" },
    @{ A = 45056.1011390162; B = "How can you help me?"; C = "It depends on your specific needs. Could you please provide more information about what you're looking for?" },
    @{ A = 45056.10121949508; B = "ransomware"; C = "We offer a free trial of our ransomware solutions so you can try them out before you buy. Would you like me to provide you with more information?" }
)

$startRow = 45
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row.A
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
}
